$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "panel_query_time" (column F) timestamps on the data sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:33:27.869946"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:27.869952"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:27.869955"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:27.869957"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:27.869959"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:27.869961"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:27.869963"
$dataSheet.Range("F9").Value = "2021-10-05 14:33:27.869965"

# --- Add the new "metadata" sheet, placed right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Copy the header / index-column formatting from the "data" sheet so the new
# sheet's styled cells (bold, centered, bordered) reuse the same style
# instead of creating brand-new ones.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Chronic granulomatous disease"
$ws.Range("C2").Value = 3159
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.1"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "2021-07-24T07:06:13.341910Z"
$ws.Range("F2").Value = "2021-10-05 14:33:27.867459"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3159/?format=json"

# Keep "data" as the active sheet, matching the original selection state.
$dataSheet.Activate()
